# Add a "pNotProductUnit" indicator column (C) before the existing "Bend"
# column (which shifts to D), fill it in with 1/0 flags, center-align the
# new column's data cells, and set the new column's width. Also move the
# active selection to A13, matching the edited workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C; existing column C ("Bend") shifts to D.
$ws.Columns("C").Insert(-4161)   # xlShiftToRight

# Header for the new column.
$ws.Range("C1").Value = "pNotProductUnit"

# Flag values for each data row (1 = non-production unit, 0 = otherwise).
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("C9").Value = 0
$ws.Range("C10").Value = 0
$ws.Range("C11").Value = 1
$ws.Range("C12").Value = 1
$ws.Range("C13").Value = 1
$ws.Range("C14").Value = 0
$ws.Range("C15").Value = 0

# Center-align the new column's data cells.
$ws.Range("C2:C15").HorizontalAlignment = -4108   # xlHAlignCenter

# New column width (closest achievable to the authored 16.140625 chars)
# and selection to match the edited workbook.
$ws.Columns("C").ColumnWidth = 15.62
$ws.Range("A13").Select()
